$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '44.450.57'
$ws.Range('E2').Value = '  +1.07%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.238.45'
$ws.Range('E3').Value = '  +0.19%  '
$ws.Range('E4').Value = '  +1.33%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '306.16'
$ws.Range('E5').Value = '  +0.29%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '94.30'
$ws.Range('E6').Value = '  -0.32%  '
$ws.Range('E7').Value = '  +0.31%  '
$ws.Range('E8').Value = '  +0.14%  '
$ws.Range('E9').Value = '  +0.30%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '34.77'
$ws.Range('E10').Value = '  +0.19%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0803'
$ws.Range('E11').Value = '  -0.32%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.21'
$ws.Range('E12').Value = '  +0.21%  '
$ws.Range('E13').Value = '  +0.06%  '
$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.834'
$ws.Range('E14').Value = '  +1.20%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.215.13'
$ws.Range('E15').Value = '  -0.92%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '13.57'
$ws.Range('E16').Value = '  -0.26%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '44.157.83'
$ws.Range('E17').Value = '  +0.62%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0₃0954'
$ws.Range('E18').Value = '  -0.36%  '
$ws.Range('E19').Value = '  +1.66%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.95'
$ws.Range('E20').Value = '  -1.26%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '65.51'
$ws.Range('E21').Value = '  +1.02%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '237.27'
$ws.Range('E22').Value = '  +0.52%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.95'
$ws.Range('E23').Value = '  +1.05%  '
$ws.Range('E24').Value = '  +1.54%  '
$ws.Range('E25').Value = '  -0.05%  '
$ws.Range('B26').Value = 'InjectiveProtocol'
$ws.Range('C26').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '38.44'
$ws.Range('E26').Value = '  +2.95%  '
$ws.Range('B27').Value = 'Toncoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.21'
$ws.Range('E27').Value = '  +2.77%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.78'
$ws.Range('E28').Value = '  -2.07%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.95'
$ws.Range('E29').Value = '  +0.49%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '19.98'
$ws.Range('E30').Value = '  +0.91%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '153.36'
$ws.Range('E31').Value = '  +0.07%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0794'
$ws.Range('E32').Value = '  -0.63%  '
$ws.Range('E33').Value = '  +2.34%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.07'
$ws.Range('E34').Value = '  -4.21%  '
$ws.Range('E35').Value = '  +4.30%  '
$ws.Range('E36').Value = '  +0.38%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.79'
$ws.Range('E37').Value = '  +2.11%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '15.13'
$ws.Range('E38').Value = '  +1.94%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.40'
$ws.Range('E39').Value = '  +1.66%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.76'
$ws.Range('E40').Value = '  -1.04%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0299'
$ws.Range('E41').Value = '  -0.34%  '
$ws.Range('E42').Value = '  +0.23%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.792.55'
$ws.Range('E43').Value = '  +3.57%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.192'
$ws.Range('E44').Value = '  +2.69%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '78.83'
$ws.Range('E45').Value = '  -8.25%  '
$ws.Range('B46').Value = 'Stacks'
$ws.Range('C46').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.62'
$ws.Range('E46').Value = '  +8.03%  '
$ws.Range('B47').Value = 'ordi'
$ws.Range('C47').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '70.35'
$ws.Range('E47').Value = '  +2.38%  '
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '98.72'
$ws.Range('E48').Value = '  -0.85%  '
$ws.Range('B49').Value = 'THORChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '4.90'
$ws.Range('E49').Value = '  +0.02%  '
$ws.Range('E50').Value = '  +0.39%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '54.50'
$ws.Range('E51').Value = '  +0.73%  '
